$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was updated from
# 45179 (2023-09-10) to 45180 (2023-09-11) for rows 2 through 27.
for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45179) {
        $cell.Value = 45180
    }
}
